$d = $word.ActiveDocument

$old = "Linear Regression (which has 1 line of best fit to identify outlier) and Higher-Dimensional Linear Regression (which the line of best fit becomes hyperplane, instead of looking at distance from a line, we look at the distance from this hyperplane)."
$new = "Logistic Regression is supervised ML algorithm used to classification. It predicts the probability hence its output lies in between 0 and 1. Random Forest it can be used for both classification and regression kind of problem. "

$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
